$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

$t.Cell(1,1).Range.Text = "40 x 10" + $vt + "  1    0" + $vt + "  ----" + $vt + "4|    |" + $vt + "0|    |"
$t.Cell(1,2).Range.Text = "26 x 87" + $vt + "  8    7" + $vt + "  ----" + $vt + "2|    |" + $vt + "6|    |"
$t.Cell(1,3).Range.Text = "95 x 44" + $vt + "  4    4" + $vt + "  ----" + $vt + "9|    |" + $vt + "5|    |"
$t.Cell(2,1).Range.Text = "40 x 80" + $vt + "  8    0" + $vt + "  ----" + $vt + "4|    |" + $vt + "0|    |"
$t.Cell(2,2).Range.Text = "94 x 53" + $vt + "  5    3" + $vt + "  ----" + $vt + "9|    |" + $vt + "4|    |"
$t.Cell(2,3).Range.Text = "71 x 15" + $vt + "  1    5" + $vt + "  ----" + $vt + "7|    |" + $vt + "1|    |"
$t.Cell(3,1).Range.Text = "73 x 48" + $vt + "  4    8" + $vt + "  ----" + $vt + "7|    |" + $vt + "3|    |"
$t.Cell(3,2).Range.Text = "80 x 32" + $vt + "  3    2" + $vt + "  ----" + $vt + "8|    |" + $vt + "0|    |"
$t.Cell(3,3).Range.Text = "14 x 90" + $vt + "  9    0" + $vt + "  ----" + $vt + "1|    |" + $vt + "4|    |"
$t.Cell(4,1).Range.Text = "95 x 58" + $vt + "  5    8" + $vt + "  ----" + $vt + "9|    |" + $vt + "5|    |"
$t.Cell(4,2).Range.Text = "16 x 71" + $vt + "  7    1" + $vt + "  ----" + $vt + "1|    |" + $vt + "6|    |"
$t.Cell(4,3).Range.Text = "14 x 44" + $vt + "  4    4" + $vt + "  ----" + $vt + "1|    |" + $vt + "4|    |"
$t.Cell(5,1).Range.Text = "79 x 65" + $vt + "  6    5" + $vt + "  ----" + $vt + "7|    |" + $vt + "9|    |"
$t.Cell(5,2).Range.Text = "50 x 15" + $vt + "  1    5" + $vt + "  ----" + $vt + "5|    |" + $vt + "0|    |"
$t.Cell(5,3).Range.Text = "35 x 27" + $vt + "  2    7" + $vt + "  ----" + $vt + "3|    |" + $vt + "5|    |"

Write-Host "done"
